$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a new bullet ("numId 2") right before the "Trang “Kiến trúc một
#    tầng”" item (numId 1), just after the "Sửa lại khung hình tầng trệt
#    khi chưa chọn." item.
# ---------------------------------------------------------------------------
$anchorIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Contains("khung hình tầng trệt khi chưa chọn")) {
        $anchorIdx = $i
    }
}
if ($anchorIdx -eq 0) {
    throw "Anchor paragraph 'khung hình tầng trệt khi chưa chọn' not found."
}

$anchorPara = $d.Paragraphs($anchorIdx)
$r = $anchorPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)

$newPara = $d.Paragraphs($anchorIdx + 1)
$newPara.Range.Text = "Sửa lại thanh định hướng (Breadcrumbs) cho đẹp hơn."

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it used to sit in front of
#    "Text canh " and now belongs in front of "Button canh phải nếu nằm
#    dưới textbox" (same run, right after the run properties).
# ---------------------------------------------------------------------------

# 2a) Stamp it onto the "Button canh phải..." run (that run is the entire
#     paragraph, so replacing the whole paragraph range keeps everything in
#     place).
$findRng = $d.Content
$found1 = $findRng.Find.Execute("Button canh phải nếu nằm dưới textbox", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
if (-not $found1) {
    throw "'Button canh phải nếu nằm dưới textbox' not found."
}
$targetRng = $d.Range($findRng.Start, $findRng.End)

$buttonXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p><w:r>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' +
  '<w:lastRenderedPageBreak/>' +
  '<w:t>Button canh phải nếu nằm dưới textbox</w:t>' +
  '</w:r></w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$targetRng.InsertXML($buttonXml)

# 2b) Remove it from the "Text canh trái nếu nằm bên trái textbox" paragraph.
#     That paragraph holds three runs, so the whole paragraph range is
#     rebuilt (in order) without the page-break marker on the first run.
$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute("Text canh ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
if (-not $found2) {
    throw "'Text canh ' not found."
}
$para2 = $findRng2.Paragraphs(1)
$para2Rng = $d.Range($para2.Range.Start, $para2.Range.End - 1)

$textCanhXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' +
  '<w:t xml:space="preserve">Text canh </w:t></w:r>' +
  '<w:r w:rsidR="00637930"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' +
  '<w:t xml:space="preserve">trái </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' +
  '<w:t>nếu nằm bên trái textbox</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$para2Rng.InsertXML($textCanhXml)
